$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: new columns I and J, with values and the same formatting as the other header cells
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

# Data values for columns I and J, rows 2-11
$dataI = @(8, 7, 8, 8, 8, 8, 8, 5, 8, 5)
$dataJ = @(9, 8, 8, 9, 9, 9, 8, 6, 8, 5)

for ($i = 0; $i -lt 10; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 9).Value = $dataI[$i]
    $ws.Cells.Item($row, 10).Value = $dataJ[$i]
}
